$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# row 1, col 1: "49÷7=7, 0" -> "23÷7=3, 2"
$tbl.Cell(1, 1).Range.Text = "23÷7=3, 2"
# row 1, col 2: "90÷5=18, 0" -> "42÷7=6, 0"
$tbl.Cell(1, 2).Range.Text = "42÷7=6, 0"
# row 1, col 3: "38÷4=9, 2" -> "22÷6=3, 4"
$tbl.Cell(1, 3).Range.Text = "22÷6=3, 4"
# row 1, col 4: "54÷6=9, 0" -> "10÷7=1, 3"
$tbl.Cell(1, 4).Range.Text = "10÷7=1, 3"
# row 1, col 5: "29÷3=9, 2" -> "63÷9=7, 0"
$tbl.Cell(1, 5).Range.Text = "63÷9=7, 0"
# row 5, col 1: "39÷5=7, 4" -> "88÷2=44, 0"
$tbl.Cell(5, 1).Range.Text = "88÷2=44, 0"
# row 5, col 2: "56÷7=8, 0" -> "60÷8=7, 4"
$tbl.Cell(5, 2).Range.Text = "60÷8=7, 4"
# row 5, col 3: "55÷7=7, 6" -> "28÷2=14, 0"
$tbl.Cell(5, 3).Range.Text = "28÷2=14, 0"
# row 5, col 4: "41÷9=4, 5" -> "13÷2=6, 1"
$tbl.Cell(5, 4).Range.Text = "13÷2=6, 1"
# row 5, col 5: "39÷3=13, 0" -> "74÷5=14, 4"
$tbl.Cell(5, 5).Range.Text = "74÷5=14, 4"
# row 9, col 1: "38÷5=7, 3" -> "26÷2=13, 0"
$tbl.Cell(9, 1).Range.Text = "26÷2=13, 0"
# row 9, col 2: "68÷9=7, 5" -> "34÷4=8, 2"
$tbl.Cell(9, 2).Range.Text = "34÷4=8, 2"
# row 9, col 3: "61÷4=15, 1" -> "41÷3=13, 2"
$tbl.Cell(9, 3).Range.Text = "41÷3=13, 2"
# row 9, col 4: "59÷5=11, 4" -> "94÷7=13, 3"
$tbl.Cell(9, 4).Range.Text = "94÷7=13, 3"
# row 9, col 5: "14÷2=7, 0" -> "27÷5=5, 2"
$tbl.Cell(9, 5).Range.Text = "27÷5=5, 2"
# row 13, col 1: "29÷5=5, 4" -> "19÷2=9, 1"
$tbl.Cell(13, 1).Range.Text = "19÷2=9, 1"
# row 13, col 2: "35÷5=7, 0" -> "89÷8=11, 1"
$tbl.Cell(13, 2).Range.Text = "89÷8=11, 1"
# row 13, col 3: "70÷3=23, 1" -> "58÷7=8, 2"
$tbl.Cell(13, 3).Range.Text = "58÷7=8, 2"
# row 13, col 4: "42÷3=14, 0" -> "70÷9=7, 7"
$tbl.Cell(13, 4).Range.Text = "70÷9=7, 7"
# row 13, col 5: "28÷2=14, 0" -> "85÷2=42, 1"
$tbl.Cell(13, 5).Range.Text = "85÷2=42, 1"
# row 17, col 1: "17÷5=3, 2" -> "99÷3=33, 0"
$tbl.Cell(17, 1).Range.Text = "99÷3=33, 0"
# row 17, col 2: "20÷2=10, 0" -> "56÷3=18, 2"
$tbl.Cell(17, 2).Range.Text = "56÷3=18, 2"
# row 17, col 3: "22÷3=7, 1" -> "19÷2=9, 1"
$tbl.Cell(17, 3).Range.Text = "19÷2=9, 1"
# row 17, col 4: "84÷6=14, 0" -> "43÷8=5, 3"
$tbl.Cell(17, 4).Range.Text = "43÷8=5, 3"
# row 17, col 5: "69÷3=23, 0" -> "36÷2=18, 0"
$tbl.Cell(17, 5).Range.Text = "36÷2=18, 0"
